# Users sheet sync: add "token" / "active_session" columns, refresh
# lastLogin timestamps + a couple of other fields (see commit "Feature/my
# order sync (#6)").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new columns before column K ("contactNumber"), which
#     pushes the existing K..O ("contactNumber".."pincode") out to M..Q.
$ws.Range("K1:L1").EntireColumn.Insert()

# --- New header cells for the inserted columns.
$ws.Range("K1").Value = "token"
$ws.Range("L1").Value = "active_session"

# Helper: write a literal TRUE/FALSE-looking (or otherwise "smart typed")
# string into a cell without Excel coercing it to a Boolean/Number. Build
# it as a formula result first, then freeze it back down to a plain value
# in place, the same way pasting "Values Only" over a formula does.
function Set-TextValue($addr, $text) {
    $cell = $ws.Range($addr)
    $escaped = $text.Replace('"', '""')
    $cell.Formula = '=""&"' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}

# --- Row 2 (superadmin_1): refreshed lastLogin + new token/active_session.
$ws.Range("J2").Value = "2025-07-23T18:15:46.122Z"
$ws.Range("K2").Value = "4c68c025-900e-4e21-902d-52f570a1d23f"
Set-TextValue "L2" "FALSE"

# --- Row 3 (Test Admin): refreshed lastLogin + new token/active_session.
$ws.Range("J3").Value = "2025-07-23T18:09:26.038Z"
$ws.Range("K3").Value = "33fdd77e-3b6b-40a3-b464-98a44403ab7c"
Set-TextValue "L3" "FALSE"

# --- Row 4 (Test Vendor): new token/active_session only.
$ws.Range("K4").Value = "a1933e0c-0a0d-4308-90a9-6a3c24dbabf0"
Set-TextValue "L4" "FALSE"

# --- Row 5 (Adarsh): refreshed lastLogin + new token/active_session.
$ws.Range("J5").Value = "2025-07-23T20:58:18.933Z"
$ws.Range("K5").Value = "084cbcfa-1f0e-49dd-b08a-463bd32f39e3"
Set-TextValue "L5" "TRUE"

# --- Row 6 (Mohit): rotated password hash, refreshed updatedAt/lastLogin,
#     new token/active_session.
$ws.Range("E6").Value = '$2a$12$RJ68RUA5TiMFSl/KCRI8j.3fQSiDivosCFmYtZQEpdXfaq6Ut2r8u'
$ws.Range("I6").Value = "2025-07-23T18:16:22.943Z"
$ws.Range("J6").Value = "2025-07-23T21:07:08.422Z"
$ws.Range("K6").Value = "883381e9-e8be-4c16-b4ed-52d51079a9f7"
Set-TextValue "L6" "TRUE"
